$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-01 Tuesday" "2024-10-02 Wednesday"

Replace-Text "493×3=" "166×4="
Replace-Text "209×3=" "188×3="
Replace-Text "994×3=" "972×3="
Replace-Text "945×7=" "236×8="
Replace-Text "809×2=" "727×8="
Replace-Text "690×9=" "949×6="
Replace-Text "757×3=" "142×5="
Replace-Text "533×5=" "705×9="
Replace-Text "443×2=" "187×7="
Replace-Text "684×6=" "420×7="
Replace-Text "829×4=" "751×3="
Replace-Text "129×2=" "462×6="
Replace-Text "817×5=" "907×4="
Replace-Text "926×9=" "222×3="
Replace-Text "607×4=" "707×4="
Replace-Text "614×4=" "273×5="
Replace-Text "478×5=" "205×9="
Replace-Text "475×2=" "976×5="
Replace-Text "519×7=" "540×6="
Replace-Text "665×4=" "268×3="
Replace-Text "837×9=" "390×9="
Replace-Text "842×5=" "524×6="
Replace-Text "960×4=" "220×7="
Replace-Text "278×2=" "847×8="
Replace-Text "952×8=" "729×5="
